$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column J: "Box Scores" style lookup/order column ---------------
$ws.Range("J5").Value  = "Spurs"
$ws.Range("J7").Value  = "Pistons"
$ws.Range("J8").Value  = "Cavaliers"
$ws.Range("J9").Value  = "Warriors"
$ws.Range("J10").Value = "Clippers"
$ws.Range("J12").Value = "Knicks"
$ws.Range("J15").Value = "Celtics"
$ws.Range("J16").Value = "Grizzlies"
$ws.Range("J22").Value = "Suns"
$ws.Range("J23").Value = "Heat"
$ws.Range("J26").Value = "Pacers"
$ws.Range("J27").Value = "Kings"
$ws.Range("J31").Value = "Bucks"

# Running count of filled team slots in column J (mirrors the other COUNTA
# helper cells already present in row 33 for columns B-F).
$ws.Range("J33").Formula = "=COUNTA(J2:J32)"

# New mode label beneath the other "Mode N" header cells in row 34.
$ws.Range("J34").Value = "Mode 6"

# Recalculate so the COUNTA formula above carries a fresh cached value.
$excel.Calculate()

# --- Hide the helper columns now that column J is the "live" view -------
$ws.Columns.Item(2).Width = 75.732421875
$ws.Columns.Item(3).Width = 51.73828125
$ws.Columns.Item(4).Width = 51.73828125
$ws.Columns.Item(5).Width = 51.73828125
$ws.Columns.Item(6).Width = 69.744140625
$ws.Columns.Item(7).Width = 51.73828125
$ws.Columns.Item(8).Width = 51.73828125

$ws.Columns.Item(2).Hidden = $true
$ws.Columns.Item(3).Hidden = $true
$ws.Columns.Item(4).Hidden = $true
$ws.Columns.Item(5).Hidden = $true
$ws.Columns.Item(6).Hidden = $true
$ws.Columns.Item(7).Hidden = $true
$ws.Columns.Item(8).Hidden = $true

# Column I becomes a zero-width hidden spacer column.
$ws.Columns.Item(9).ColumnWidth = -0.8333333333333334
$ws.Columns.Item(9).Hidden = $true

# --- Move the active selection to K5, matching the author's last click --
$ws.Range("K5").Select()
